# Fruta / hortaliza, semanal
# Inserts two new weekly price rows (variety "Kakamas") at the top of the
# data block (rows 113-114), pushing the existing rows 113-124 down to
# 115-126 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 113; everything currently at 113.. shifts to 115..
$ws.Range("A113:T114").EntireRow.Insert()

# Row 113: Kakamas / Especial
$ws.Range("A113").Value = 2
$ws.Range("B113").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C113").Value = "Coquimbo"
$ws.Range("D113").Value = 44637
$ws.Range("E113").Value = 4
$ws.Range("F113").Value = "Fruta"
$ws.Range("G113").Value = 100103
$ws.Range("H113").Value = "Frutos de hueso (carozo)"
$ws.Range("I113").Value = 100103004
$ws.Range("J113").Value = "Durazno"
$ws.Range("K113").Value = "Kakamas"
$ws.Range("L113").Value = "Especial"
$ws.Range("M113").Value = 20
$ws.Range("N113").Value = 460000
$ws.Range("O113").Value = 470000
$ws.Range("P113").Value = 465000
$ws.Range("Q113").Value = "`$/bins (400 kilos)"
$ws.Range("R113").Value = "Región de O'Higgins"
$ws.Range("S113").Value = 1162
$ws.Range("T113").Value = 400

# Row 114: Kakamas / Primera
$ws.Range("A114").Value = 2
$ws.Range("B114").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C114").Value = "Coquimbo"
$ws.Range("D114").Value = 44637
$ws.Range("E114").Value = 4
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100103
$ws.Range("H114").Value = "Frutos de hueso (carozo)"
$ws.Range("I114").Value = 100103004
$ws.Range("J114").Value = "Durazno"
$ws.Range("K114").Value = "Kakamas"
$ws.Range("L114").Value = "Primera"
$ws.Range("M114").Value = 20
$ws.Range("N114").Value = 410000
$ws.Range("O114").Value = 420000
$ws.Range("P114").Value = 415000
$ws.Range("Q114").Value = "`$/bins (400 kilos)"
$ws.Range("R114").Value = "Región de O'Higgins"
$ws.Range("S114").Value = 1038
$ws.Range("T114").Value = 400
